# Update Mappings 22 Ontologies
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Stash the pristine "Hyperlink" cell format (s="2") in a scratch cell far
# away from the used range, before touching any hyperlinks - Hyperlinks.Add
# re-applies the Hyperlink style itself and (re)creates a slightly different
# style entry, so we restore the original look afterwards from this stash.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial($xlPasteFormats)

# --- New header cell F1 (OBI_DEF) - copy formatting from E1 ---
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial($xlPasteFormats)
$ws.Range("F1").Value = "OBI_DEF"

# Fix existing header labels (metadata4Ing_* -> metadata4ing_*)
$ws.Range("B1").Value = "metadata4ing_IRI"
$ws.Range("C1").Value = "metadata4ing_DESC"

# --- Formatting for the two brand-new rows (6 and 7) ---
# Column A (index numbers)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial($xlPasteFormats)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A7").PasteSpecial($xlPasteFormats)

# Column B/D (hyperlink styled)
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B6").PasteSpecial($xlPasteFormats)
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteFormats)
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D7").PasteSpecial($xlPasteFormats)

# --- Row 2: BFO_0000015 (Process) - add new OBI_DEF text ---
$ws.Range("F2").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(`"Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'`", 'en')]"

# --- Row 3: BFO_0000017 (realizable entity) - add new OBI_DEF text ---
$ws.Range("F3").Value = "['To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]', 'To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])']"

# --- Row 4: new "Group" row (inserted before Organization) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "http://xmlns.com/foaf/0.1/Group"
$ws.Range("C4").Value = "{'label': 'Group', 'prefLabel': 'Group', 'altLabel': None, 'name': 'Group'}"
$ws.Range("D4").Value = "http://purl.obolibrary.org/obo/OBI_0302900"
$ws.Range("E4").Value = "{'label': 'Group'}"
$ws.Range("F4").Value = "[]"

# --- Row 5: Organization (was row4, values updated, hyperlink location removed) ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "http://xmlns.com/foaf/0.1/Organization"
$ws.Range("C5").Value = "{'label': 'Organization', 'prefLabel': 'Organisation', 'altLabel': None, 'name': 'Organization'}"
$ws.Range("D5").Value = "http://purl.obolibrary.org/obo/OBI_0000245"
$ws.Range("E5").Value = "{'label': 'Organization'}"
$ws.Range("F5").Value = "[]"

# --- Row 6: new "Dataset" row ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "http://www.w3.org/ns/dcat#Dataset"
$ws.Range("C6").Value = "{'label': None, 'prefLabel': 'Datensatz', 'altLabel': None, 'name': 'Dataset'}"
$ws.Range("D6").Value = "http://purl.obolibrary.org/obo/APOLLO_SV_00000796"
$ws.Range("E6").Value = "{'label': 'Dataset'}"
$ws.Range("F6").Value = "[]"

# --- Row 7: Role (was row5, description dict updated, new OBI_DEF text) ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "http://www.w3.org/ns/prov#Role"
$ws.Range("C7").Value = "{'label': None, 'prefLabel': 'Role', 'altLabel': None, 'name': 'Role'}"
$ws.Range("D7").Value = "http://purl.obolibrary.org/obo/BFO_0000023"
$ws.Range("E7").Value = "{'label': 'Role', 'prefLabel': 'Role'}"
$ws.Range("F7").Value = "['B is a role means: b is a realizable entity and b exists because there is some single bearer that is in some special physical, social, or institutional set of circumstances in which this bearer does not have to be and b is not such that, if it ceases to exist, then the physical make-up of the bearer is thereby changed. [BFO]']"

# --- Hyperlinks: clear all existing and re-add in the final order so the ---
# --- relationship ids line up the way Excel would renumber them.        ---
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "http://purl.obolibrary.org/obo/BFO_0000015") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "http://purl.obolibrary.org/obo/BFO_0000015") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "http://purl.obolibrary.org/obo/BFO_0000017") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "http://purl.obolibrary.org/obo/BFO_0000017") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "http://xmlns.com/foaf/0.1/Group") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "http://purl.obolibrary.org/obo/OBI_0302900") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "http://xmlns.com/foaf/0.1/Organization") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "http://purl.obolibrary.org/obo/OBI_0000245") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "http://www.w3.org/ns/dcat", "Dataset") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "http://purl.obolibrary.org/obo/APOLLO_SV_00000796") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "http://www.w3.org/ns/prov", "Role") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "http://purl.obolibrary.org/obo/BFO_0000023") | Out-Null

# Restore the original Hyperlink cell look (undoes the extra style variant
# that Hyperlinks.Add creates) on every linked cell.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B2").PasteSpecial($xlPasteFormats)
$ws.Range("D2").PasteSpecial($xlPasteFormats)
$ws.Range("B3").PasteSpecial($xlPasteFormats)
$ws.Range("D3").PasteSpecial($xlPasteFormats)
$ws.Range("B4").PasteSpecial($xlPasteFormats)
$ws.Range("D4").PasteSpecial($xlPasteFormats)
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("D5").PasteSpecial($xlPasteFormats)
$ws.Range("B6").PasteSpecial($xlPasteFormats)
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("B7").PasteSpecial($xlPasteFormats)
$ws.Range("D7").PasteSpecial($xlPasteFormats)
$ws.Range("Z1").Clear()
